$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Cells.Item(2, 4)
$cell.Value = "'52.030.97"
$cell.Style = "Normal"

$cell = $ws.Cells.Item(3, 4)
$cell.Value = "'2.832.65"
$cell.Style = "Normal"
$ws.Cells.Item(3, 5).Value = "  +2.72%  "

$cell = $ws.Cells.Item(4, 4)
$cell.Value = "'1.00"
$cell.Style = "Normal"
$ws.Cells.Item(4, 5).Value = "  -0.02%  "

$cell = $ws.Cells.Item(5, 4)
$cell.Value = "'355.73"
$cell.Style = "Normal"
$ws.Cells.Item(5, 5).Value = "  +7.05%  "

$cell = $ws.Cells.Item(6, 4)
$cell.Value = "'114.29"
$cell.Style = "Normal"
$ws.Cells.Item(6, 5).Value = "  -1.82%  "

$cell = $ws.Cells.Item(7, 4)
$cell.Value = "'0.554"
$cell.Style = "Normal"
$ws.Cells.Item(7, 5).Value = "  +2.83%  "

$ws.Cells.Item(8, 5).Value = "  -0.05%  "

$cell = $ws.Cells.Item(9, 4)
$cell.Value = "'0.603"
$cell.Style = "Normal"
$ws.Cells.Item(9, 5).Value = "  +4.76%  "

$cell = $ws.Cells.Item(10, 4)
$cell.Value = "'41.73"
$cell.Style = "Normal"
$ws.Cells.Item(10, 5).Value = "  +0.02%  "

$cell = $ws.Cells.Item(11, 4)
$cell.Value = "'0.0855"
$cell.Style = "Normal"
$ws.Cells.Item(11, 5).Value = "  -0.26%  "

$cell = $ws.Cells.Item(12, 4)
$cell.Value = "'20.07"
$cell.Style = "Normal"
$ws.Cells.Item(12, 5).Value = "  -0.61%  "

$ws.Cells.Item(13, 5).Value = "  +1.57%  "

$cell = $ws.Cells.Item(14, 4)
$cell.Value = "'7.77"
$cell.Style = "Normal"
$ws.Cells.Item(14, 5).Value = "  +1.59%  "

$cell = $ws.Cells.Item(15, 4)
$cell.Value = "'3.264.49"
$cell.Style = "Normal"
$ws.Cells.Item(15, 5).Value = "  +2.29%  "

$cell = $ws.Cells.Item(16, 4)
$cell.Value = "'2.826.00"
$cell.Style = "Normal"
$ws.Cells.Item(16, 5).Value = "  +2.10%  "

$cell = $ws.Cells.Item(17, 4)
$cell.Value = "'0.900"
$cell.Style = "Normal"
$ws.Cells.Item(17, 5).Value = "  +1.47%  "

$cell = $ws.Cells.Item(18, 4)
$cell.Value = "'51.902.96"
$cell.Style = "Normal"
$ws.Cells.Item(18, 5).Value = "  +0.71%  "

$cell = $ws.Cells.Item(19, 4)
$cell.Value = "'7.42"
$cell.Style = "Normal"
$ws.Cells.Item(19, 5).Value = "  +8.10%  "

$ws.Cells.Item(20, 5).Value = "  -1.76%  "

$cell = $ws.Cells.Item(21, 4)
$cell.Value = "'13.59"
$cell.Style = "Normal"
$ws.Cells.Item(21, 5).Value = "  +0.85%  "

$cell = $ws.Cells.Item(22, 4)
$cell.Value = "'0.0000100"
$cell.Style = "Normal"
$ws.Cells.Item(22, 5).Value = "  +2.39%  "

$cell = $ws.Cells.Item(23, 4)
$cell.Value = "'270.57"
$cell.Style = "Normal"
$ws.Cells.Item(23, 5).Value = "  -2.75%  "

$cell = $ws.Cells.Item(24, 4)
$cell.Value = "'69.97"
$cell.Style = "Normal"
$ws.Cells.Item(24, 5).Value = "  +0.34%  "

$ws.Cells.Item(25, 5).Value = "  +5.34%  "

$cell = $ws.Cells.Item(26, 4)
$cell.Value = "'26.82"
$cell.Style = "Normal"
$ws.Cells.Item(26, 5).Value = "  -0.02%  "

$ws.Cells.Item(27, 5).Value = "  -0.04%  "

$cell = $ws.Cells.Item(28, 4)
$cell.Value = "'10.35"
$cell.Style = "Normal"
$ws.Cells.Item(28, 5).Value = "  +1.40%  "

$ws.Cells.Item(29, 5).Value = "  +1.54%  "

$ws.Cells.Item(30, 5).Value = "  -0.77%  "

$cell = $ws.Cells.Item(31, 4)
$cell.Value = "'0.0457"
$cell.Style = "Normal"
$ws.Cells.Item(31, 5).Value = "  +30.16%  "

$cell = $ws.Cells.Item(32, 4)
$cell.Value = "'50.72"
$cell.Style = "Normal"
$ws.Cells.Item(32, 5).Value = "  +1.59%  "

$cell = $ws.Cells.Item(33, 4)
$cell.Value = "'33.99"
$cell.Style = "Normal"
$ws.Cells.Item(33, 5).Value = "  -3.24%  "

$cell = $ws.Cells.Item(34, 4)
$cell.Value = "'5.84"
$cell.Style = "Normal"
$ws.Cells.Item(34, 5).Value = "  +5.08%  "

$cell = $ws.Cells.Item(35, 4)
$cell.Value = "'0.0834"
$cell.Style = "Normal"
$ws.Cells.Item(35, 5).Value = "  +1.09%  "

$ws.Cells.Item(36, 5).Value = "  -0.10%  "

$cell = $ws.Cells.Item(37, 4)
$cell.Value = "'2.09"
$cell.Style = "Normal"
$ws.Cells.Item(37, 5).Value = "  +0.28%  "

$ws.Cells.Item(38, 2).Value = "RenderToken"
$ws.Cells.Item(38, 3).Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$cell = $ws.Cells.Item(38, 4)
$cell.Value = "'4.93"
$cell.Style = "Normal"
$ws.Cells.Item(38, 5).Value = "  -1.10%  "

$ws.Cells.Item(39, 2).Value = "LidoDAOToken"
$ws.Cells.Item(39, 3).Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$cell = $ws.Cells.Item(39, 4)
$cell.Value = "'3.23"
$cell.Style = "Normal"
$ws.Cells.Item(39, 5).Value = "  -0.11%  "

$cell = $ws.Cells.Item(40, 4)
$cell.Value = "'18.22"
$cell.Style = "Normal"
$ws.Cells.Item(40, 5).Value = "  -4.37%  "

$cell = $ws.Cells.Item(41, 4)
$cell.Value = "'23.70"
$cell.Style = "Normal"
$ws.Cells.Item(41, 5).Value = "  +2.58%  "

$cell = $ws.Cells.Item(42, 4)
$cell.Value = "'2.58"
$cell.Style = "Normal"
$ws.Cells.Item(42, 5).Value = "  +6.10%  "

$ws.Cells.Item(43, 5).Value = "  +1.39%  "

$cell = $ws.Cells.Item(44, 4)
$cell.Value = "'126.44"
$cell.Style = "Normal"
$ws.Cells.Item(44, 5).Value = "  -0.17%  "

$ws.Cells.Item(45, 5).Value = "  +0.37%  "

$cell = $ws.Cells.Item(46, 4)
$cell.Value = "'3.38"
$cell.Style = "Normal"
$ws.Cells.Item(46, 5).Value = "  +1.77%  "

$cell = $ws.Cells.Item(47, 4)
$cell.Value = "'2.088.95"
$cell.Style = "Normal"
$ws.Cells.Item(47, 5).Value = "  -0.06%  "

$cell = $ws.Cells.Item(49, 4)
$cell.Value = "'5.74"
$cell.Style = "Normal"
$ws.Cells.Item(49, 5).Value = "  +3.55%  "

$cell = $ws.Cells.Item(50, 4)
$cell.Value = "'0.953"
$cell.Style = "Normal"
$ws.Cells.Item(50, 5).Value = "  +9.30%  "

$cell = $ws.Cells.Item(51, 4)
$cell.Value = "'60.87"
$cell.Style = "Normal"
$ws.Cells.Item(51, 5).Value = "  +1.68%  "
